$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.699.72'
$ws.Range('D3').Value = '3.118.24'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.32'
$ws.Range('E5').Value = '  +9.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '631.83'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('E7').Value = '  +1.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.362'
$ws.Range('E8').Value = '  -2.84%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').Value = '3.114.63'
$ws.Range('E10').Value = '  -1.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.717'
$ws.Range('E11').Value = '  -4.15%  '
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.63'
$ws.Range('E13').Value = '  +5.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000246'
$ws.Range('E14').Value = '  -1.64%  '
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('D16').Value = '90.598.40'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Value = '3.691.93'
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').Value = '3.115.12'
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.80'
$ws.Range('E19').Value = '  +1.78%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.16'
$ws.Range('E20').Value = '  -0.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000209'
$ws.Range('E21').Value = '  -4.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '442.27'
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('E23').Value = '  +6.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.95'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.91'
$ws.Range('E25').Value = '  -3.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.58'
$ws.Range('E26').Value = '  +1.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '88.27'
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('D28').Value = '3.305.46'
$ws.Range('E28').Value = '  +0.24%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.54'
$ws.Range('E30').Value = '  +3.78%  '
$ws.Range('E31').Value = '  -2.76%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.195'
$ws.Range('E32').Value = '  +24.88%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.30'
$ws.Range('E33').Value = '  +3.78%  '
$ws.Range('B34').Value = 'dogwifhat'
$ws.Range('C34').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.82'
$ws.Range('E34').Value = '  +2.44%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.889'
$ws.Range('E35').Value = '  -1.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '508.99'
$ws.Range('E36').Value = '  -3.95%  '
$ws.Range('E37').Value = '  +3.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.14'
$ws.Range('E38').Value = '  +1.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.92'
$ws.Range('E39').Value = '  +1.39%  '
$ws.Range('E40').Value = '  -1.91%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.411'
$ws.Range('E41').Value = '  +1.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.16'
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0846'
$ws.Range('E44').Value = '  +3.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.27'
$ws.Range('E45').Value = '  +48.60%  '
$ws.Range('E46').Value = '  -1.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '151.31'
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.686'
$ws.Range('E48').Value = '  +5.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '45.17'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.48'
$ws.Range('E51').Value = '  +2.23%  '
